$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "sortOrder" column (E) is being removed. Before removing it, its values
# become the new "id" values (column B), replacing the old id numbering.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $sortOrderValue = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 2).Value = $sortOrderValue
}

# Now delete the sortOrder column (E) entirely, shifting F:N left to E:M
$ws.Columns.Item(5).Delete()

# Update the active cell selection like the committed workbook
$ws.Range("F27").Select()
